# Fruta / hortaliza, semanal
# The weekly refresh rotates the per-row data (date, variety, quality,
# volume, min/max/avg price, origin region, $/Kg) among rows 2-9 while
# leaving the market/region/product descriptor columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each data row (row -> column values), matching the
# post-edit state described by the diff.
$rows = @{
    2 = @{ D = 44355; K = "Mankaki";  L = "Segunda"; M = 270; N = 20000; O = 21000; P = 20500; R = "Región Metropolitana";  S = 1139 }
    3 = @{ D = 44301; K = "Hachiya";  L = "Segunda"; M = 250; N = 20000; O = 21000; P = 20500; R = "Región de O'Higgins";   S = 1139 }
    4 = @{ D = 44313; K = "Mankaki";  L = "Primera"; M = 270; N = 21000; O = 22000; P = 21500; R = "Región de O'Higgins";   S = 1194 }
    5 = @{ D = 45043; K = "Fuyu";     L = "Primera"; M = 300; N = 25000; O = 26000; P = 25500; R = "Región de O'Higgins";   S = 1417 }
    6 = @{ D = 44305; K = "Mankaki";  L = "Segunda"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins";   S = 1361 }
    7 = @{ D = 44342; K = "Mankaki";  L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins";   S = 1361 }
    8 = @{ D = 45071; K = "Fuyu";     L = "Segunda"; M = 110; N = 23000; O = 24000; P = 23455; R = "Región Metropolitana";  S = 1303 }
    9 = @{ D = 44699; K = "Mankaki";  L = "Primera"; M = 250; N = 29000; O = 30000; P = 29500; R = "Región de O'Higgins";   S = 1639 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
}
